# Fruta / hortaliza, semanal
# Insert 3 new weekly price rows for "Palta" (avocado) before the current
# row 919, shifting the existing data (rows 919-996) down to rows 922-999.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at 919..921 (pushes old 919.. downward)
$ws.Range("A919:T921").EntireRow.Insert()

# Populate the 3 new rows with this week's data (constant columns A,B,C,E,F,G,H,I,J
# mirror the rest of the sheet; D is the new reporting date 2022-09-22 = serial 44826)

# Row 919: Hass - Especial
$ws.Range("A919").Value = 8
$ws.Range("B919").Value = "Terminal La Palmera de La Serena"
$ws.Range("C919").Value = "Coquimbo"
$ws.Range("D919").Value = 44826
$ws.Range("E919").Value = 4
$ws.Range("F919").Value = "Fruta"
$ws.Range("G919").Value = 100106
$ws.Range("H919").Value = "Oleaginosos"
$ws.Range("I919").Value = 100106002
$ws.Range("J919").Value = "Palta"
$ws.Range("K919").Value = "Hass"
$ws.Range("L919").Value = "Especial"
$ws.Range("M919").Value = 480
$ws.Range("N919").Value = 21000
$ws.Range("O919").Value = 22000
$ws.Range("P919").Value = 21500
$ws.Range("Q919").Value = "$/bandeja 10 kilos"
$ws.Range("R919").Value = "Perú"
$ws.Range("S919").Value = 2150
$ws.Range("T919").Value = 10

# Row 920: Hass - Primera
$ws.Range("A920").Value = 8
$ws.Range("B920").Value = "Terminal La Palmera de La Serena"
$ws.Range("C920").Value = "Coquimbo"
$ws.Range("D920").Value = 44826
$ws.Range("E920").Value = 4
$ws.Range("F920").Value = "Fruta"
$ws.Range("G920").Value = 100106
$ws.Range("H920").Value = "Oleaginosos"
$ws.Range("I920").Value = 100106002
$ws.Range("J920").Value = "Palta"
$ws.Range("K920").Value = "Hass"
$ws.Range("L920").Value = "Primera"
$ws.Range("M920").Value = 400
$ws.Range("N920").Value = 19000
$ws.Range("O920").Value = 20000
$ws.Range("P920").Value = 19500
$ws.Range("Q920").Value = "$/bandeja 10 kilos"
$ws.Range("R920").Value = "Perú"
$ws.Range("S920").Value = 1950
$ws.Range("T920").Value = 10

# Row 921: Hass - Segunda
$ws.Range("A921").Value = 8
$ws.Range("B921").Value = "Terminal La Palmera de La Serena"
$ws.Range("C921").Value = "Coquimbo"
$ws.Range("D921").Value = 44826
$ws.Range("E921").Value = 4
$ws.Range("F921").Value = "Fruta"
$ws.Range("G921").Value = 100106
$ws.Range("H921").Value = "Oleaginosos"
$ws.Range("I921").Value = 100106002
$ws.Range("J921").Value = "Palta"
$ws.Range("K921").Value = "Hass"
$ws.Range("L921").Value = "Segunda"
$ws.Range("M921").Value = 300
$ws.Range("N921").Value = 16000
$ws.Range("O921").Value = 17000
$ws.Range("P921").Value = 16500
$ws.Range("Q921").Value = "$/bandeja 10 kilos"
$ws.Range("R921").Value = "Perú"
$ws.Range("S921").Value = 1650
$ws.Range("T921").Value = 10
